# Trade #4 closed at 2026-02-17 08:07:48 - unknown UNKNOWN +0.000%
#
# This adds a new closed trade (#4) for the MarketMaking strategy to the
# "All Trades" and "MarketMaking" sheets, and updates the aggregated
# metrics on the "Summary" and "Strategy Status" sheets accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.85   # Current Capital
$summary.Range("B4").Value = -0.15     # Total P&L $
$summary.Range("B5").Value = -0.75     # Total P&L %
$summary.Range("B6").Value = 4         # Total Trades
$summary.Range("B8").Value = 2         # Losing Trades
$summary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.84999999999999  # Capital
$status.Range("D4").Value = 4                  # Trades
$status.Range("E4").Value = -0.15              # P&L $
$status.Range("F4").Value = -0.15              # P&L %
$status.Range("G4").Value = 25                 # Win Rate %

# ---------------------------------------------------------------------
# Helper to append the new trade row (#4) to a trade-log sheet that uses
# the same column layout ("All Trades" and "MarketMaking" sheets).
# ---------------------------------------------------------------------
function Add-TradeRow4($ws) {
    $row = 5

    $ws.Cells.Item($row, 1).Value = 4
    # Date / Time columns hold text that looks like a date/time, so force
    # text interpretation with a leading quote-prefix (same as typing
    # '2026-02-17 into Excel) to avoid automatic date conversion.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 3).Value = "'08:07:41"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.79
    $ws.Cells.Item($row, 7).Value = 0.6899999999999999
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -12.6582
    $ws.Cells.Item($row, 10).Value = -0.1
    $ws.Cells.Item($row, 11).Value = 99.84999999999999
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

Add-TradeRow4 $wb.Worksheets.Item("All Trades")
Add-TradeRow4 $wb.Worksheets.Item("MarketMaking")
